$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column Z (26) for "Logical Processor"
$ws.Columns.Item(26).Insert()

# Header for new column Z
$ws.Cells.Item(1, 26).Value = "Logical Processor"
$ws.Cells.Item(1, 26).Style = $ws.Cells.Item(1, 25).Style

# Data value for new column Z
$ws.Cells.Item(2, 26).Value = 16
$ws.Cells.Item(2, 26).Style = $ws.Cells.Item(2, 25).Style

# Set column width to match diff (bestFit width ~16.14)
$ws.Columns.Item(26).ColumnWidth = 16.140625

# Update selection / view
$ws.Range("I9").Select()
